$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: renamed extr1 -> line7, C/D updated, E unchanged (TRUE)
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9: renamed extr2 -> line8, C updated, D unchanged, E flips FALSE->TRUE
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Row 10: renamed extr3 -> extr1, C/D updated, E unchanged (TRUE)
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11: renamed extr4 -> extr2, C/D updated, E unchanged (TRUE)
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12: renamed extr5 -> extr3, C updated, D unchanged, E flips FALSE->TRUE
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

# Row 13: renamed extr6 -> extr4, C unchanged, D updated, E unchanged (FALSE)
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

# Row 14: renamed extr7 -> extr5, C/D updated, E unchanged (TRUE)
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# Row 15: renamed extr8 -> extr6, C/D updated, E unchanged (FALSE)
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# Row 16: new row, extr7
$ws.Range("A16").Value = 14
$ws.Range("A16").Style = $ws.Range("A15").Style
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# Row 17: new row, extr8
$ws.Range("A17").Value = 15
$ws.Range("A17").Style = $ws.Range("A15").Style
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
